# This script applies the edits described by the commit:
#  "add new Panama species list, revise code to check it, and revised temp check files"
#
# Concretely, for this workbook (tocheck/temptnrsprob.xlsx) that means:
#  1. The row for "Brugmansia candida" (row 4) was removed from the TNRS
#     results, shifting every subsequent row up by one.
#  2. A new column "dateTNRS" (AW) was appended, populated with the date the
#     TNRS check was run (2025-09-01) for every data row, formatted as
#     yyyy-mm-dd, and the column was widened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Brugmansia candida" row (row 4); rows below shift up.
$ws.Rows("4").Delete()

# 2. Add the new "dateTNRS" column (AW) with a header and a date value for
#    every remaining data row (rows 2-18 after the deletion above).
$ws.Range("AW1").Value = "dateTNRS"

$lastRow = $ws.Range("A1").End(4).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("AW$r")
    $cell.Value = 45901
    $cell.NumberFormat = "yyyy-mm-dd"
}

# 3. Widen the new column to fit the date values.
$ws.Columns("AW").ColumnWidth = 19.86
